$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 70 (shifts existing rows 70..133 down to 71..134)
$ws.Rows.Item(70).Insert()

# Populate the newly inserted row 70 with this week's data (same static
# market/category metadata as the surrounding rows, new date and prices)
$ws.Cells.Item(70, 1).Value = 8
$ws.Cells.Item(70, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(70, 3).Value = "Coquimbo"
$ws.Cells.Item(70, 4).Value = 44651
$ws.Cells.Item(70, 4).NumberFormat = $ws.Cells.Item(71, 4).NumberFormat
$ws.Cells.Item(70, 5).Value = 4
$ws.Cells.Item(70, 6).Value = 100112040
$ws.Cells.Item(70, 7).Value = "Cilantro"
$ws.Cells.Item(70, 8).Value = "Sin especificar"
$ws.Cells.Item(70, 9).Value = "Primera"
$ws.Cells.Item(70, 10).Value = 2400
$ws.Cells.Item(70, 11).Value = 2000
$ws.Cells.Item(70, 12).Value = 2500
$ws.Cells.Item(70, 13).Value = 2250
$ws.Cells.Item(70, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(70, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(70, 16).Value = 1500
$ws.Cells.Item(70, 17).Value = 1.5
$ws.Cells.Item(70, 18).Value = "Hortaliza"
